$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 209
$ws.Range("J6").Value = 266.5
$ws.Range("L6").Value = 799.5
$ws.Range("N6").Value = -1023.5

$ws.Range("H8").Value = 2
$ws.Range("I8").Value = 2
$ws.Range("K8").Value = 6
$ws.Range("M8").Value = 133

$ws.Range("H33").Value = 13141.883
$ws.Range("I33").Value = 16462.309
$ws.Range("K33").Value = 16462.309
$ws.Range("M33").Value = -16233.309

$ws.Range("H42").Value = 219.83333
$ws.Range("I42").Value = 48.285713
$ws.Range("J42").Value = 460
$ws.Range("K42").Value = 144.857139
$ws.Range("L42").Value = 1380
$ws.Range("M42").Value = 85.14286099999998
$ws.Range("N42").Value = -1840

$ws.Range("H51").Value = 4732.933
$ws.Range("I51").Value = 3666.6667
$ws.Range("J51").Value = 4999.5
$ws.Range("K51").Value = 3666.6667
$ws.Range("L51").Value = 4999.5
$ws.Range("M51").Value = -3182.6667
$ws.Range("N51").Value = -5967.5

$ws.Range("H64").Value = 6984.6562
$ws.Range("I64").Value = 5809.381
$ws.Range("J64").Value = 9228.362999999999
$ws.Range("K64").Value = 5809.381
$ws.Range("L64").Value = 9228.362999999999
$ws.Range("M64").Value = -5561.381
$ws.Range("N64").Value = -9724.362999999999

$ws.Range("H67").Value = 6984.6562
$ws.Range("I67").Value = 5809.381
$ws.Range("J67").Value = 9228.362999999999
$ws.Range("K67").Value = 5809.381
$ws.Range("L67").Value = 9228.362999999999
$ws.Range("M67").Value = -4951.381
$ws.Range("N67").Value = -10944.363

$ws.Range("H138").Value = 34488416
$ws.Range("J138").Value = 58832360
$ws.Range("L138").Value = 176497080
$ws.Range("N138").Value = -176507360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1063.52
$ws.Range("I32").Value = 909.7527
$ws.Range("J32").Value = 3106.4285
$ws.Range("K32").Value = 909.7527
$ws.Range("L32").Value = 3106.4285
$ws.Range("M32").Value = -622.7527
$ws.Range("N32").Value = -3680.4285

$ws.Range("H74").Value = 100113410
$ws.Range("I74").Value = 166854340
$ws.Range("K74").Value = 166854340
$ws.Range("M74").Value = -166853466

$ws.Range("H77").Value = 100113410
$ws.Range("I77").Value = 166854340
$ws.Range("K77").Value = 834271700
$ws.Range("M77").Value = -834267332

$ws.Range("H97").Value = 1945.2632
$ws.Range("I97").Value = 1747.7142
$ws.Range("K97").Value = 1747.7142
$ws.Range("M97").Value = -1251.7142

$ws.Range("H122").Value = 2152.5
$ws.Range("I122").Value = 1430.1428
$ws.Range("K122").Value = 4290.428400000001
$ws.Range("M122").Value = -1840.428400000001

$ws.Range("H132").Value = 16170817
$ws.Range("I132").Value = 3378.4902
$ws.Range("K132").Value = 10135.4706
$ws.Range("M132").Value = -7605.470600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1564.3103
$ws.Range("I20").Value = 1658.0526
$ws.Range("J20").Value = 1386.2
$ws.Range("K20").Value = 1658.0526
$ws.Range("L20").Value = 1386.2
$ws.Range("M20").Value = -1411.0526
$ws.Range("N20").Value = -1880.2

$ws.Range("H134").Value = 2502378
$ws.Range("I134").Value = 2633503
$ws.Range("K134").Value = 7900509
$ws.Range("M134").Value = -7897974

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 388.45456
$ws.Range("I22").Value = 382.3
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 382.3
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = -32.30000000000001
$ws.Range("N22").Value = -1150

$ws.Range("H31").Value = 32261336
$ws.Range("I31").Value = 2303.6
$ws.Range("K31").Value = 2303.6
$ws.Range("M31").Value = -2008.6

$ws.Range("H34").Value = 32261336
$ws.Range("I34").Value = 2303.6
$ws.Range("K34").Value = 2303.6
$ws.Range("M34").Value = -2101.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 865.8
$ws.Range("I119").Value = 865.8
$ws.Range("K119").Value = 2597.4
$ws.Range("M119").Value = 2240.6

$ws.Range("H131").Value = 1477.4286
$ws.Range("J131").Value = 1696.84
$ws.Range("L131").Value = 5090.52
$ws.Range("N131").Value = -15170.52

$ws.Range("H134").Value = 1028.7576
$ws.Range("I134").Value = 1028.7576
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3086.2728
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 1983.7272
$ws.Range("N134").ClearContents()

$ws.Range("H138").Value = 2122.3845
$ws.Range("I138").Value = 1721.1111
$ws.Range("K138").Value = 5163.3333
$ws.Range("M138").Value = -23.33330000000024

$ws.Range("H140").Value = 3366
$ws.Range("I140").Value = 1850
$ws.Range("K140").Value = 5550
$ws.Range("M140").Value = -370

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 17298.8
$ws.Range("I21").Value = 17623.5
$ws.Range("K21").Value = 17623.5
$ws.Range("M21").Value = -17450.5

$ws.Range("H30").Value = 17298.8
$ws.Range("I30").Value = 17623.5
$ws.Range("K30").Value = 17623.5
$ws.Range("M30").Value = -17518.5

$ws.Range("H70").Value = 7204
$ws.Range("I70").Value = 7204
$ws.Range("K70").Value = 7204
$ws.Range("M70").Value = -6934

$ws.Range("H73").Value = 7204
$ws.Range("I73").Value = 7204
$ws.Range("K73").Value = 7204
$ws.Range("M73").Value = -6268

$ws.Range("H97").Value = 2602.3684
$ws.Range("J97").Value = 3502.5
$ws.Range("L97").Value = 3502.5
$ws.Range("N97").Value = -4494.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2675.8572
$ws.Range("I22").Value = 1840.2222
$ws.Range("J22").Value = 4180
$ws.Range("K22").Value = 1840.2222
$ws.Range("L22").Value = 4180
$ws.Range("M22").Value = -1545.2222
$ws.Range("N22").Value = -4770

$ws.Range("H27").Value = 2675.8572
$ws.Range("I27").Value = 1840.2222
$ws.Range("J27").Value = 4180
$ws.Range("K27").Value = 1840.2222
$ws.Range("L27").Value = 4180
$ws.Range("M27").Value = -1733.2222
$ws.Range("N27").Value = -4394

$ws.Range("H46").Value = 1830.4348
$ws.Range("J46").Value = 3427.7778
$ws.Range("L46").Value = 3427.7778
$ws.Range("N46").Value = -3803.7778

$ws.Range("H55").Value = 655.73334
$ws.Range("J55").Value = 803.36365
$ws.Range("L55").Value = 803.36365
$ws.Range("N55").Value = -1149.36365

$ws.Range("H61").Value = 5788.7617
$ws.Range("I61").Value = 4118.933
$ws.Range("K61").Value = 4118.933
$ws.Range("M61").Value = -3916.933

$ws.Range("H113").Value = 5788.7617
$ws.Range("I113").Value = 4118.933
$ws.Range("K113").Value = 4118.933
$ws.Range("M113").Value = -1948.933

$ws.Range("H132").Value = 40007940
$ws.Range("I132").Value = 4025.139
$ws.Range("J132").Value = 142875150
$ws.Range("K132").Value = 12075.417
$ws.Range("L132").Value = 428625450
$ws.Range("M132").Value = -9545.417000000001
$ws.Range("N132").Value = -428630510

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 56500
$ws.Range("I68").Value = 56500
$ws.Range("K68").Value = 56500
$ws.Range("M68").Value = -55689

$ws.Range("H71").Value = 56500
$ws.Range("I71").Value = 56500
$ws.Range("K71").Value = 169500
$ws.Range("M71").Value = -165444

$ws.Range("H96").Value = 3299.6
$ws.Range("I96").Value = 1666.3334
$ws.Range("J96").Value = 5749.5
$ws.Range("K96").Value = 1666.3334
$ws.Range("L96").Value = 5749.5
$ws.Range("M96").Value = -293.3334
$ws.Range("N96").Value = -8495.5

$ws.Range("H132").Value = 2574.4138
$ws.Range("I132").Value = 2354.56
$ws.Range("K132").Value = 7063.68
$ws.Range("M132").Value = -4533.68

$ws.Range("H136").Value = 1965.2727
$ws.Range("I136").Value = 1021.72
$ws.Range("K136").Value = 3065.16
$ws.Range("M136").Value = -515.1599999999999

$ws.Range("H141").Value = 67282.42999999999
$ws.Range("J141").Value = 81399.39999999999
$ws.Range("L141").Value = 81399.39999999999
$ws.Range("N141").Value = -91759.39999999999
